$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7-25 down to 8-26
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new Guayaba price record
$ws.Range("A7").Value = 10
$ws.Range("B7").Value = "Vega Modelo de Temuco"
$ws.Range("C7").Value = "La Araucanía"
$ws.Range("D7").Value = 45054
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100108
$ws.Range("H7").Value = "Tropicales y subtropicales"
$ws.Range("I7").Value = 100108001
$ws.Range("J7").Value = "Guayaba"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 25
$ws.Range("N7").Value = 2500
$ws.Range("O7").Value = 2500
$ws.Range("P7").Value = 2500
$ws.Range("Q7").Value = "$/kilo"
$ws.Range("R7").Value = "Región de Arica y Parinacota"
$ws.Range("S7").Value = 2500
$ws.Range("T7").Value = 1
